$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three "ECs -> ..." sender rows (old rows 2-4) are being dropped; the
# three "MuSCs -> ..." sender rows (old rows 5-7) become the new rows 2-4,
# with refreshed TPM-derived numbers. Delete the old rows 5-7 first (so the
# addresses of rows 2-4 below stay put), then delete the old ECs rows 2-4,
# which shifts the former rows 5-7 up into 2-4.
$ws.Rows("5:7").Delete()
$ws.Rows("2:4").Delete()

# Row 2: MuSCs / Areg / Erbb3 / ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Areg"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014328
$ws.Range("H2").Value = 0.042984
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09551033333333332
$ws.Range("N2").Value = 0.286531
$ws.Range("O2").Value = 0.0198020999427218
$ws.Range("P2").Value = 0.0198020999427218
$ws.Range("Q2").Value = 0.001368472056
$ws.Range("R2").Value = 0.012316248504
$ws.Range("S2").Value = 0.0198020999427218
$ws.Range("T2").Value = 0.0198020999427218

# Row 3: MuSCs / Areg / Erbb3 / FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Areg"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014328
$ws.Range("H3").Value = 0.042984
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3460976666666666
$ws.Range("N3").Value = 1.038293
$ws.Range("O3").Value = 0.07175622098770619
$ws.Range("P3").Value = 0.07175622098770619
$ws.Range("Q3").Value = 0.004958887368
$ws.Range("R3").Value = 0.044629986312
$ws.Range("S3").Value = 0.07175622098770619
$ws.Range("T3").Value = 0.07175622098770619

# Row 4: MuSCs / Areg / Erbb3 / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Areg"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.014328
$ws.Range("H4").Value = 0.042984
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.381634666666667
$ws.Range("N4").Value = 13.144904
$ws.Range("O4").Value = 0.9084416790695721
$ws.Range("P4").Value = 0.9084416790695721
$ws.Range("Q4").Value = 0.06278006150400001
$ws.Range("R4").Value = 0.565020553536
$ws.Range("S4").Value = 0.9084416790695721
$ws.Range("T4").Value = 0.9084416790695721
